# Update the crypto price/volume table with the latest scraped values
# (GitHub Actions refresh). Column D ("Price") holds text like "27.900.28"
# that Excel would otherwise auto-parse as a number, so numeric-looking
# values are written with a leading quote to force text entry, matching
# the existing t="inlineStr"/text-typed cells in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.900.28'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').Value = '1.646.62'
$ws.Range('E3').Value = '  +1.78%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''213.54'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '''23.56'
$ws.Range('E8').Value = '  +3.88%  '
$ws.Range('D9').Value = '''0.266'
$ws.Range('E9').Value = '  +1.44%  '
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('E11').Value = '  -1.64%  '
$ws.Range('D12').Value = '1.880.45'
$ws.Range('E12').Value = '  +1.79%  '
$ws.Range('D13').Value = '1.650.40'
$ws.Range('E13').Value = '  +1.95%  '
$ws.Range('E14').Value = '  +1.56%  '
$ws.Range('E15').Value = '  +2.24%  '
$ws.Range('D16').Value = '''65.66'
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('D17').Value = '27.907.33'
$ws.Range('E17').Value = '  +1.40%  '
$ws.Range('D18').Value = '''231.84'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('D19').Value = '''7.66'
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('D20').Value = '0.0₃0724'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').Value = '''10.71'
$ws.Range('E22').Value = '  +5.96%  '
$ws.Range('E23').Value = '  +2.35%  '
$ws.Range('E24').Value = '  +2.83%  '
$ws.Range('D25').Value = '''152.10'
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('E26').Value = '  +1.91%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''15.73'
$ws.Range('E27').Value = '  +1.27%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '''0.112'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = '''1.20'
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('E31').Value = '  +0.83%  '
$ws.Range('E32').Value = '  +2.78%  '
$ws.Range('D33').Value = '1.454.70'
$ws.Range('E33').Value = '  +0.95%  '
$ws.Range('E34').Value = '  +2.23%  '
$ws.Range('E35').Value = '  +1.87%  '
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('D37').Value = '''0.889'
$ws.Range('E37').Value = '  +3.45%  '
$ws.Range('E38').Value = '  +1.01%  '
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('D40').Value = '''0.922'
$ws.Range('E40').Value = '  -1.16%  '
$ws.Range('D41').Value = '''69.33'
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('E42').Value = '  +2.27%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('E45').Value = '  +1.13%  '
$ws.Range('E47').Value = '  +5.61%  '
$ws.Range('D48').Value = '1.789.18'
$ws.Range('E48').Value = '  +1.68%  '
$ws.Range('D49').Value = '''88.95'
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('D50').Value = '0.0₆0106'
$ws.Range('E50').Value = '  +0.75%  '
$ws.Range('E51').Value = '  +1.11%  '
